$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "sex": same number of rows (3 data rows), values/order change so that
# row2=Non-binary, row3=Man, row4=Woman.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("sex")
$sexData = @(
    @("Non-binary", 1,  0,    -0.5413533834586466, 0,                   -0.9927536231884058, [double]"2.377935733887305e-10", -0.02671138269588797),
    @("Man",        76, 0.4473684210526316, -0.093984962406015, 1,      0.007246376811594235, 0.03308232690219727,  0.006370943968515726),
    @("Woman",      189,0.582010582010582,  0.04065719855193539, 0.9905660377358491, -0.002187585452556706, 0.02429085193415082, -0.002420530999530723)
)
$r = 2
foreach ($row in $sexData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet "glang": shrinks from 13 data rows (rows 2-14) to 10 data rows (rows
# 2-12). Delete the now-unused rows 13:14 then overwrite rows 2:12.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("glang")
$ws.Range("A13:H14").ClearContents()

$glangData = @(
    @("Croatian",   1,  0,   -0.5413533834586466, 0, -0.9927536231884058, 0.0006130258047512081, -0.02609835712893033),
    @("German",     8,  0.375, -0.1663533834586466, 0.75, -0.2427536231884058, 0.04952669376063773, 0.02281531082695619),
    @("English",    10, 0.4, -0.1413533834586466, 1, 0.007246376811594235, 0.007208463115822105, -0.01950291981785944),
    @("Italian",    16, 0.5, -0.04135338345864659, 1, 0.007246376811594235, 0.01091044529753344, -0.0158009376361481),
    @("Romanian",   2,  0.5, -0.04135338345864659, 1, 0.007246376811594235, 0.0002163129610653677, -0.02649506997261617),
    @("French",     212,0.5283018867924528, -0.01305149666619376, 1, 0.007246376811594235, 0.02901499345507749, 0.002303610521395948),
    @("Other",      4,  0.75, 0.2086466165413534, 1, 0.007246376811594235, 0.02575229467046115, -0.0009590882632203937),
    @("Arab",       2,  1,   0.4586466165413534, 1, 0.007246376811594235, [double]"2.183694983766856e-05", -0.02668954598384388),
    @("Portuguese", 9,  1,   0.4586466165413534, 1, 0.007246376811594235, 0.02203319115004268, -0.00467819178363886),
    @("Russian",    1,  1,   0.4586466165413534, 1, 0.007246376811594235, 0.008786489787561962, -0.01792489314611958),
    @("Turkish",    1,  1,   0.4586466165413534, 1, 0.007246376811594235, [double]"2.274326564742327e-07", -0.02671115550102507)
)
$r = 2
foreach ($row in $glangData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet "glang_group": 2 data rows, swap German/Non-German order & values.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("glang_group")
$glangGroupData = @(
    @("German",     8,   0.375,              -0.1663533834586466, 0.75, -0.2427536231884058,  0.04952669376063773, 0.02281531082695619),
    @("Non-German", 258, 0.5465116279069767, 0.005158244448330129, 1,    0.007246376811594235, 0.02600393143517128, -0.0007074514985102656)
)
$r = 2
foreach ($row in $glangGroupData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet "year": 6 data rows, reorder/update values (labels mostly the same
# set, with row order changed and row7 "Bmed2" staying put with new values).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("year")
$yearData = @(
    @("Mmed3", 29, 0.2413793103448276, -0.299974073113819,    1,                  0.007246376811594235, 0.009841289206476313, -0.01687009372720523),
    @("Mmed1", 40, 0.45,               -0.09135338345864658,  1,                  0.007246376811594235, 0.01272748592356357,  -0.01398389701011797),
    @("Mmed2", 40, 0.45,               -0.09135338345864658,  1,                  0.007246376811594235, 0.02181056550004147,  -0.004900817433640074),
    @("Bmed3", 49, 0.5714285714285714, 0.03007518796992481,   1,                  0.007246376811594235, 0.03513777886965821,  0.008426395935976671),
    @("Bmed1", 67, 0.6268656716417911, 0.08551228818314449,   0.975609756097561, -0.01714386709084481,  0.03254799485066859,  0.005836611916987051),
    @("Bmed2", 41, 0.7560975609756098, 0.2147441775169632,    1,                  0.007246376811594235, 0.03745955114193237,  0.01074816820825083)
)
$r = 2
foreach ($row in $yearData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}
